# Auto-generated Excel COM-interop script
# Applies the scheduled-runner price/profit refresh to the 8 Leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per the commit diff.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 30310258
$ws.Range("I64").Value = 47625948
$ws.Range("J64").Value = 7799.25
$ws.Range("K64").Value = 47625948
$ws.Range("L64").Value = 7799.25
$ws.Range("M64").Value = -47625700
$ws.Range("N64").Value = -8295.25
$ws.Range("H67").Value = 30310258
$ws.Range("I67").Value = 47625948
$ws.Range("J67").Value = 7799.25
$ws.Range("K67").Value = 47625948
$ws.Range("L67").Value = 7799.25
$ws.Range("M67").Value = -47625090
$ws.Range("N67").Value = -9515.25
$ws.Range("H76").Value = 31256122
$ws.Range("I76").Value = 6150
$ws.Range("K76").Value = 6150
$ws.Range("M76").Value = -5835
$ws.Range("H79").Value = 31256122
$ws.Range("I79").Value = 6150
$ws.Range("K79").Value = 6150
$ws.Range("M79").Value = -5058
$ws.Range("H80").Value = 37762.285
$ws.Range("I80").Value = 25192.25
$ws.Range("J80").Value = 42790.3
$ws.Range("K80").Value = 75576.75
$ws.Range("L80").Value = 128370.9
$ws.Range("M80").Value = -74578.75
$ws.Range("N80").Value = -130366.9
$ws.Range("H83").Value = 37762.285
$ws.Range("I83").Value = 25192.25
$ws.Range("J83").Value = 42790.3
$ws.Range("K83").Value = 226730.25
$ws.Range("L83").Value = 385112.7
$ws.Range("M83").Value = -221738.25
$ws.Range("N83").Value = -395096.7
$ws.Range("H132").Value = 871.0714
$ws.Range("I132").Value = 896.6667
$ws.Range("J132").Value = 180
$ws.Range("K132").Value = 2690.0001
$ws.Range("L132").Value = 540
$ws.Range("M132").Value = -160.0001000000002
$ws.Range("N132").Value = -5600
$ws.Range("H138").Value = 2045374.4
$ws.Range("J138").Value = 3454414.5
$ws.Range("L138").Value = 10363243.5
$ws.Range("N138").Value = -10373523.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H16").Value = 6680
$ws.Range("I16").Value = 1133.3334
$ws.Range("K16").Value = 1133.3334
$ws.Range("M16").Value = -846.3334
$ws.Range("H32").Value = 4006533
$ws.Range("I32").Value = 4171596.8
$ws.Range("K32").Value = 4171596.8
$ws.Range("M32").Value = -4171309.8
$ws.Range("H61").Value = 3953.84
$ws.Range("I61").Value = 2830.3333
$ws.Range("J61").Value = 4990.923
$ws.Range("K61").Value = 2830.3333
$ws.Range("L61").Value = 4990.923
$ws.Range("M61").Value = -2618.3333
$ws.Range("N61").Value = -5414.923
$ws.Range("H74").Value = 22008.66
$ws.Range("I74").Value = 25104.904
$ws.Range("J74").Value = 5753.375
$ws.Range("K74").Value = 25104.904
$ws.Range("L74").Value = 5753.375
$ws.Range("M74").Value = -24230.904
$ws.Range("N74").Value = -7501.375
$ws.Range("H77").Value = 22008.66
$ws.Range("I77").Value = 25104.904
$ws.Range("J77").Value = 5753.375
$ws.Range("K77").Value = 125524.52
$ws.Range("L77").Value = 28766.875
$ws.Range("M77").Value = -121156.52
$ws.Range("N77").Value = -37502.875
$ws.Range("H105").Value = 99000
$ws.Range("J105").Value = 99000
$ws.Range("L105").Value = 99000
$ws.Range("N105").Value = -105988
$ws.Range("H132").Value = 3558.4102
$ws.Range("I132").Value = 1447.7307
$ws.Range("J132").Value = 7779.769
$ws.Range("K132").Value = 4343.1921
$ws.Range("L132").Value = 23339.307
$ws.Range("M132").Value = -1813.1921
$ws.Range("N132").Value = -28399.307
$ws.Range("H136").Value = 3953.84
$ws.Range("I136").Value = 2830.3333
$ws.Range("J136").Value = 4990.923
$ws.Range("K136").Value = 8490.999899999999
$ws.Range("L136").Value = 14972.769
$ws.Range("M136").Value = -5940.999899999999
$ws.Range("N136").Value = -20072.769

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 30786
$ws.Range("I26").Value = 24997
$ws.Range("K26").Value = 24997
$ws.Range("M26").Value = -24705
$ws.Range("H94").Value = 1537.0769
$ws.Range("I94").Value = 692
$ws.Range("J94").Value = 5400.2856
$ws.Range("K94").Value = 692
$ws.Range("L94").Value = 5400.2856
$ws.Range("M94").Value = -241
$ws.Range("N94").Value = -6302.2856
$ws.Range("H107").Value = 40187990
$ws.Range("I107").Value = 53581104
$ws.Range("K107").Value = 53581104
$ws.Range("M107").Value = -53579184
$ws.Range("H134").Value = 7147039
$ws.Range("I134").Value = 13159667
$ws.Range("J134").Value = 7043.6875
$ws.Range("K134").Value = 39479001
$ws.Range("L134").Value = 21131.0625
$ws.Range("M134").Value = -39476466
$ws.Range("N134").Value = -26201.0625

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4620.029
$ws.Range("I31").Value = 2303.6316
$ws.Range("J31").Value = 7370.75
$ws.Range("K31").Value = 2303.6316
$ws.Range("L31").Value = 7370.75
$ws.Range("M31").Value = -2008.6316
$ws.Range("N31").Value = -7960.75
$ws.Range("H34").Value = 4620.029
$ws.Range("I34").Value = 2303.6316
$ws.Range("J34").Value = 7370.75
$ws.Range("K34").Value = 2303.6316
$ws.Range("L34").Value = 7370.75
$ws.Range("M34").Value = -2101.6316
$ws.Range("N34").Value = -7774.75
$ws.Range("H58").Value = 4554.484
$ws.Range("I58").Value = 1810.6154
$ws.Range("J58").Value = 6536.1665
$ws.Range("K58").Value = 1810.6154
$ws.Range("L58").Value = 6536.1665
$ws.Range("M58").Value = -1607.6154
$ws.Range("N58").Value = -6942.1665
$ws.Range("H132").Value = 3326.6584
$ws.Range("I132").Value = 2474.16
$ws.Range("J132").Value = 4658.6875
$ws.Range("K132").Value = 7422.48
$ws.Range("L132").Value = 13976.0625
$ws.Range("M132").Value = -4892.48
$ws.Range("N132").Value = -19036.0625
$ws.Range("H134").Value = 3866.25
$ws.Range("I134").Value = 2097.65
$ws.Range("K134").Value = 6292.950000000001
$ws.Range("M134").Value = -3757.950000000001
$ws.Range("H136").Value = 4554.484
$ws.Range("I136").Value = 1810.6154
$ws.Range("J136").Value = 6536.1665
$ws.Range("K136").Value = 5431.8462
$ws.Range("L136").Value = 19608.4995
$ws.Range("M136").Value = -2881.8462
$ws.Range("N136").Value = -24708.4995

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 682
$ws.Range("J12").Value = 35.727272
$ws.Range("L12").Value = 107.181816
$ws.Range("N12").Value = -453.181816
$ws.Range("H14").Value = 166666660
$ws.Range("I14").Value = 166666660
$ws.Range("K14").Value = 499999980
$ws.Range("M14").Value = -499999807

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("N21").Value = ""
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("N30").Value = ""
$ws.Range("H92").Value = 18949
$ws.Range("J92").Value = 18949
$ws.Range("L92").Value = 18949
$ws.Range("N92").Value = -22693
$ws.Range("H132").Value = 1969.28
$ws.Range("I132").Value = 1650.75
$ws.Range("J132").Value = 3243.4
$ws.Range("K132").Value = 4952.25
$ws.Range("L132").Value = 9730.200000000001
$ws.Range("M132").Value = -2422.25
$ws.Range("N132").Value = -14790.2
$ws.Range("H136").Value = 29255.545
$ws.Range("J136").Value = 29090.111
$ws.Range("L136").Value = 87270.333
$ws.Range("N136").Value = -92370.333

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5055.56
$ws.Range("J7").Value = 5925.6665
$ws.Range("L7").Value = 5925.6665
$ws.Range("N7").Value = -6149.6665
$ws.Range("H106").Value = 20729.727
$ws.Range("J106").Value = 20729.727
$ws.Range("L106").Value = 20729.727
$ws.Range("N106").Value = -23253.727
$ws.Range("H122").Value = 3795.5
$ws.Range("I122").Value = 3795.5
$ws.Range("K122").Value = 11386.5
$ws.Range("M122").Value = -8936.5
$ws.Range("H126").Value = 5055.56
$ws.Range("J126").Value = 5925.6665
$ws.Range("L126").Value = 17776.9995
$ws.Range("N126").Value = -22716.9995
$ws.Range("H132").Value = 7252705.5
$ws.Range("I132").Value = 12197501
$ws.Range("K132").Value = 36592503
$ws.Range("M132").Value = -36589973
$ws.Range("H136").Value = 8434.472
$ws.Range("I136").Value = 2418.138
$ws.Range("J136").Value = 15704.208
$ws.Range("K136").Value = 7254.414
$ws.Range("L136").Value = 47112.624
$ws.Range("M136").Value = -4704.414
$ws.Range("N136").Value = -52212.624

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 762.0833
$ws.Range("I100").Value = 294.14285
$ws.Range("K100").Value = 588.2857
$ws.Range("M100").Value = -47.28570000000002
$ws.Range("H104").Value = 54175
$ws.Range("J104").Value = 54175
$ws.Range("L104").Value = 54175
$ws.Range("N104").Value = -61163
$ws.Range("H122").Value = 10288639
$ws.Range("I122").Value = 13623447
$ws.Range("K122").Value = 40870341
$ws.Range("M122").Value = -40867891
$ws.Range("H132").Value = 5297.2856
$ws.Range("I132").Value = 5169.091
$ws.Range("K132").Value = 15507.273
$ws.Range("M132").Value = -12977.273

